# Monthly Billable Hours Report — add a "Target(hrs)" column to the report
# and strip the header row off the clipboard export so the pasted TSV lines
# up with the template's data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Rename the "Target" header to "Target(hrs)" to match the other hrs columns.
$ws.Range("F1").Value = "Target(hrs)"

# The placeholder 160-hrs/month target values were just sample data for the
# chart; clear them out so the column is blank until the user fills in real
# target hours.
$ws.Range("F2:F13").ClearContents()

# Widen column F now that the header text is longer.
$ws.Columns.Item(6).ColumnWidth = 9.75

# The exported TSV (for pasting into the sheet) shouldn't include its own
# header row, since row 1 here already has headers - skip the first line
# before it hits the clipboard.
$ws.Range("A15").Value = 'Get-TogglMonthlyHoursReport | ConvertTo-Csv -NoTypeInformation -Delimiter "`t" | select -Skip 1 | clip'
